$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI LR-pairs results updated following Dr Hou's advice:
# Ligand-expressing-cells count changed from 1 to 3 for all Sending clusters,
# the Target cluster dimension now includes both "FAPs" and "sCs" (previously
# only "sCs" was present), which expands the data from 4 rows (one per Sending
# cluster) to 8 rows (4 Sending clusters x 2 Target clusters), and all derived
# statistics (expression values, specificities, weights) were recomputed.

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Lpar3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 162.399297
$ws.Range("H2").Value = 487.197891
$ws.Range("I2").Value = 0.3910371682630009
$ws.Range("J2").Value = 0.3910371682630009
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04225933333333334
$ws.Range("N2").Value = 0.126778
$ws.Range("O2").Value = 0.02137494054695286
$ws.Range("P2").Value = 0.02137494054695286
$ws.Range("Q2").Value = 6.862886025022
$ws.Range("R2").Value = 61.765974225198
$ws.Range("S2").Value = 0.008358396223270446
$ws.Range("T2").Value = 0.008358396223270446

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Lpar3"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 162.399297
$ws.Range("H3").Value = 487.197891
$ws.Range("I3").Value = 0.3910371682630009
$ws.Range("J3").Value = 0.3910371682630009
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.934791
$ws.Range("N3").Value = 5.804373
$ws.Range("O3").Value = 0.9786250594530471
$ws.Range("P3").Value = 0.9786250594530471
$ws.Range("Q3").Value = 314.208698241927
$ws.Range("R3").Value = 2827.878284177343
$ws.Range("S3").Value = 0.3826787720397305
$ws.Range("T3").Value = 0.3826787720397305

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Lpar3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 65.41736466666667
$ws.Range("H4").Value = 196.252094
$ws.Range("I4").Value = 0.1575168212364948
$ws.Range("J4").Value = 0.1575168212364948
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04225933333333334
$ws.Range("N4").Value = 0.126778
$ws.Range("O4").Value = 0.02137494054695286
$ws.Range("P4").Value = 0.02137494054695286
$ws.Range("Q4").Value = 2.764494219236889
$ws.Range("R4").Value = 24.880447973132
$ws.Range("S4").Value = 0.003366912689075078
$ws.Range("T4").Value = 0.003366912689075078

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Lpar3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 65.41736466666667
$ws.Range("H5").Value = 196.252094
$ws.Range("I5").Value = 0.1575168212364948
$ws.Range("J5").Value = 0.1575168212364948
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.934791
$ws.Range("N5").Value = 5.804373
$ws.Range("O5").Value = 0.9786250594530471
$ws.Range("P5").Value = 0.9786250594530471
$ws.Range("Q5").Value = 126.5689284007847
$ws.Range("R5").Value = 1139.120355607062
$ws.Range("S5").Value = 0.1541499085474198
$ws.Range("T5").Value = 0.1541499085474198

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Lpar3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 126.3069433333333
$ws.Range("H6").Value = 378.92083
$ws.Range("I6").Value = 0.3041313008456065
$ws.Range("J6").Value = 0.3041313008456065
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04225933333333334
$ws.Range("N6").Value = 0.126778
$ws.Range("O6").Value = 0.02137494054695286
$ws.Range("P6").Value = 0.02137494054695286
$ws.Range("Q6").Value = 5.337647220637778
$ws.Range("R6").Value = 48.03882498574001
$ws.Range("S6").Value = 0.006500788474042272
$ws.Range("T6").Value = 0.006500788474042274

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Lpar3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 126.3069433333333
$ws.Range("H7").Value = 378.92083
$ws.Range("I7").Value = 0.3041313008456065
$ws.Range("J7").Value = 0.3041313008456065
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.934791
$ws.Range("N7").Value = 5.804373
$ws.Range("O7").Value = 0.9786250594530471
$ws.Range("P7").Value = 0.9786250594530471
$ws.Range("Q7").Value = 244.3775371988433
$ws.Range("R7").Value = 2199.39783478959
$ws.Range("S7").Value = 0.2976305123715642
$ws.Range("T7").Value = 0.2976305123715642

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Lpar3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 61.180387
$ws.Range("H8").Value = 183.541161
$ws.Range("I8").Value = 0.1473147096548978
$ws.Range("J8").Value = 0.1473147096548978
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04225933333333334
$ws.Range("N8").Value = 0.126778
$ws.Range("O8").Value = 0.02137494054695286
$ws.Range("P8").Value = 0.02137494054695286
$ws.Range("Q8").Value = 2.585442367695333
$ws.Range("R8").Value = 23.268981309258
$ws.Range("S8").Value = 0.003148843160565063
$ws.Range("T8").Value = 0.003148843160565064

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Lpar3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 61.180387
$ws.Range("H9").Value = 183.541161
$ws.Range("I9").Value = 0.1473147096548978
$ws.Range("J9").Value = 0.1473147096548978
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.934791
$ws.Range("N9").Value = 5.804373
$ws.Range("O9").Value = 0.9786250594530471
$ws.Range("P9").Value = 0.9786250594530471
$ws.Range("Q9").Value = 118.371262144117
$ws.Range("R9").Value = 1065.341359297053
$ws.Range("S9").Value = 0.1441658664943327
$ws.Range("T9").Value = 0.1441658664943327

